$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 2023-09-01 (45170) to 2023-09-05 (45174) for rows 2-27
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45174
}
